$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B41").Value = "ROANOKE"
$ws.Range("B50").Value = "STLOUIS"
